# Chinh sua ke hoach
# Mark a handful of "Ngoc Nhi" / "Phuoc Toan" task cells as done ("x"),
# and refresh the saved window scroll position / active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E22").Value = "x"
$ws.Range("D23").Value = "x"
$ws.Range("D25").Value = "x"
$ws.Range("E26").Value = "x"
$ws.Range("D27").Value = "x"
$ws.Range("E28").Value = "x"
$ws.Range("D29").Value = "x"
$ws.Range("E30").Value = "x"

# Update the view: scroll the window so row 7 is the top-left visible row,
# and leave the active selection on E30.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E30").Select()
